# Apply scheduled market-price refresh to the Spriggan Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per leve row,
# sheet by sheet (one worksheet per crafting class).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 3
$ws.Range("H3").Value = 11379.6
$ws.Range("J3").Value = 11379.6
$ws.Range("L3").Value = 11379.6
$ws.Range("N3").Value = -11607.6
# row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
# row 17
$ws.Range("H17").Value = 334998
$ws.Range("J17").Value = 334998
$ws.Range("L17").Value = 1004994
$ws.Range("N17").Value = -1005330
# row 76
$ws.Range("H76").Value = 8535
$ws.Range("I76").Value = 7912.2856
$ws.Range("J76").Value = 9624.75
$ws.Range("K76").Value = 7912.2856
$ws.Range("L76").Value = 9624.75
$ws.Range("M76").Value = -7597.2856
$ws.Range("N76").Value = -10254.75
# row 79
$ws.Range("H79").Value = 8535
$ws.Range("I79").Value = 7912.2856
$ws.Range("J79").Value = 9624.75
$ws.Range("K79").Value = 7912.2856
$ws.Range("L79").Value = 9624.75
$ws.Range("M79").Value = -6820.2856
$ws.Range("N79").Value = -11808.75
# row 102
$ws.Range("H102").Value = 11379.6
$ws.Range("J102").Value = 11379.6
$ws.Range("L102").Value = 11379.6
$ws.Range("N102").Value = -17869.6
# row 113
$ws.Range("H113").Value = 1794.6471
$ws.Range("I113").Value = 1709.7142
$ws.Range("J113").Value = 2191
$ws.Range("K113").Value = 1709.7142
$ws.Range("L113").Value = 2191
$ws.Range("M113").Value = 1544.2858
$ws.Range("N113").Value = -8699
# row 131
$ws.Range("H131").Value = 1090
$ws.Range("I131").Value = 948
$ws.Range("K131").Value = 2844
$ws.Range("M131").Value = 2196

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 5631.7236
$ws.Range("I32").Value = 3366.17
$ws.Range("K32").Value = 3366.17
$ws.Range("M32").Value = -3079.17
# row 104
$ws.Range("H104").Value = 20176.334
$ws.Range("J104").Value = 25162.5
$ws.Range("L104").Value = 25162.5
$ws.Range("N104").Value = -32150.5
# row 122
$ws.Range("H122").Value = 1800.8096
$ws.Range("I122").Value = 2025.1538
$ws.Range("J122").Value = 1436.25
$ws.Range("K122").Value = 6075.4614
$ws.Range("L122").Value = 4308.75
$ws.Range("M122").Value = -3625.4614
$ws.Range("N122").Value = -9208.75
# row 132
$ws.Range("H132").Value = 4354100.5
$ws.Range("I132").Value = 5006371
$ws.Range("J132").Value = 5629.6665
$ws.Range("K132").Value = 15019113
$ws.Range("L132").Value = 16888.9995
$ws.Range("M132").Value = -15016583
$ws.Range("N132").Value = -21948.9995

$ws = $wb.Worksheets.Item("CRP")
# row 17
$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5348
# row 99
$ws.Range("H99").Value = 3843.375
$ws.Range("I99").Value = 4100.5386
$ws.Range("J99").Value = 2729
$ws.Range("K99").Value = 4100.5386
$ws.Range("L99").Value = 2729
$ws.Range("M99").Value = -2602.5386
$ws.Range("N99").Value = -5725
# row 122
$ws.Range("H122").Value = 3786.5
$ws.Range("I122").Value = 3597.463
$ws.Range("J122").Value = 5062.5
$ws.Range("K122").Value = 10792.389
$ws.Range("L122").Value = 15187.5
$ws.Range("M122").Value = -8342.389000000001
$ws.Range("N122").Value = -20087.5
# row 126
$ws.Range("H126").Value = 3843.375
$ws.Range("I126").Value = 4100.5386
$ws.Range("J126").Value = 2729
$ws.Range("K126").Value = 12301.6158
$ws.Range("L126").Value = 8187
$ws.Range("M126").Value = -9831.6158
$ws.Range("N126").Value = -13127
# row 141
$ws.Range("H141").Value = 208364.53
$ws.Range("J141").Value = 263245.94
$ws.Range("L141").Value = 263245.94
$ws.Range("N141").Value = -273605.94

$ws = $wb.Worksheets.Item("CUL")
# row 8
$ws.Range("H8").Value = 326.27274
$ws.Range("I8").Value = 326.27274
$ws.Range("K8").Value = 978.81822
$ws.Range("M8").Value = -839.81822
# row 9
$ws.Range("H9").Value = 3222
$ws.Range("I9").Value = 2333
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 6999
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = -6775
$ws.Range("N9").Value = -15448
# row 58
$ws.Range("H58").Value = 500
$ws.Range("J58").Value = 500
$ws.Range("L58").Value = 1500
$ws.Range("N58").Value = -1756
# row 69
$ws.Range("H69").Value = 808
$ws.Range("I69").Value = 720.5
$ws.Range("K69").Value = 2161.5
$ws.Range("M69").Value = -1350.5
# row 72
$ws.Range("H72").Value = 808
$ws.Range("I72").Value = 720.5
$ws.Range("K72").Value = 6484.5
$ws.Range("M72").Value = -2428.5

$ws = $wb.Worksheets.Item("GSM")
# row 7
$ws.Range("H7").Value = 3350001
$ws.Range("I7").Value = 10000000
$ws.Range("K7").Value = 10000000
$ws.Range("M7").Value = -9999888
# row 8
$ws.Range("H8").Value = 3350001
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999861
# row 12
$ws.Range("H12").Value = 5469.75
$ws.Range("I12").Value = 2989.5
$ws.Range("J12").Value = 7950
$ws.Range("K12").Value = 2989.5
$ws.Range("L12").Value = 7950
$ws.Range("M12").Value = -2849.5
$ws.Range("N12").Value = -8230
# row 70
$ws.Range("H70").Value = 14288
$ws.Range("J70").Value = 13077
$ws.Range("L70").Value = 13077
$ws.Range("N70").Value = -13617
# row 73
$ws.Range("H73").Value = 14288
$ws.Range("J73").Value = 13077
$ws.Range("L73").Value = 13077
$ws.Range("N73").Value = -14949
# row 80
$ws.Range("H80").Value = 1267.3889
$ws.Range("I80").Value = 1038.375
$ws.Range("J80").Value = 1450.6
$ws.Range("K80").Value = 1038.375
$ws.Range("L80").Value = 1450.6
$ws.Range("M80").Value = -40.375
$ws.Range("N80").Value = -3446.6
# row 83
$ws.Range("H83").Value = 1267.3889
$ws.Range("I83").Value = 1038.375
$ws.Range("J83").Value = 1450.6
$ws.Range("K83").Value = 5191.875
$ws.Range("L83").Value = 7253
$ws.Range("M83").Value = -199.875
$ws.Range("N83").Value = -17237
# row 134
$ws.Range("H134").Value = 19500
$ws.Range("J134").Value = 19500
$ws.Range("L134").Value = 58500
$ws.Range("N134").Value = -63570

$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 1884.6538
$ws.Range("K61").Value = 1884.6538
$ws.Range("M61").Value = -1682.6538
# row 95
$ws.Range("H95").Value = 32200
$ws.Range("J95").Value = 32200
$ws.Range("L95").Value = 32200
$ws.Range("N95").Value = -37692
# row 100
$ws.Range("H100").Value = 14586933
$ws.Range("I100").Value = 15912109
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 15912109
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -15911568
$ws.Range("N100").Value = -11082
# row 113
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 1884.6538
$ws.Range("K113").Value = 1884.6538
$ws.Range("M113").Value = 285.3462

$ws = $wb.Worksheets.Item("WVR")
# row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# row 4
$ws.Range("H4").Value = 100006060
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100006060
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100006060
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -100006286
# row 93
$ws.Range("H93").Value = 39999
$ws.Range("J93").Value = 39999
$ws.Range("L93").Value = 39999
$ws.Range("N93").Value = -44991
# row 96
$ws.Range("H96").Value = 2572.5833
$ws.Range("J96").Value = 2781.6155
$ws.Range("L96").Value = 2781.6155
$ws.Range("N96").Value = -5527.6155
